# Applies the "Sat Oct 14 21:30:28 UTC 2023" cryptos refresh: updated prices/
# 1h-volume deltas for every row, plus three coins that swapped rank order
# (rows 35/36, 40/41, 49/50 exchange Coin/Link/Price/Volume together).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values are plain text in this sheet (prices use "." as a thousands
# separator, e.g. "26.942.92"); force text format first so COM's type
# auto-detection does not coerce them to numbers and drop trailing zeros
# (e.g. "62.10" -> 62.1).
$updates = @{
    "D2" = "26.942.92"
    "E2" = "  -0.61%  "
    "D3" = "1.562.36"
    "E3" = "  -0.42%  "
    "D4" = "1.01"
    "E4" = "  +0.18%  "
    "D5" = "207.54"
    "E5" = "  -0.27%  "
    "E6" = "  -0.13%  "
    "E7" = "  +0.03%  "
    "D8" = "22.08"
    "E8" = "  -0.08%  "
    "E9" = "  -0.77%  "
    "E10" = "  +1.87%  "
    "E11" = "  -0.31%  "
    "D12" = "1.785.68"
    "E12" = "  -0.30%  "
    "D13" = "1.563.55"
    "E13" = "  -0.41%  "
    "E14" = "  -0.16%  "
    "E15" = "  -0.43%  "
    "D16" = "62.10"
    "E16" = "  +0.00%  "
    "D17" = "26.946.26"
    "E17" = "  -0.58%  "
    "E18" = "  +1.01%  "
    "D19" = "216.86"
    "E19" = "  -0.98%  "
    "E20" = "  -0.01%  "
    "E21" = "  +0.04%  "
    "E22" = "  +0.55%  "
    "E23" = "  -0.87%  "
    "E24" = "  -1.16%  "
    "D25" = "152.50"
    "E25" = "  -1.20%  "
    "E26" = "  -0.38%  "
    "D27" = "15.06"
    "E27" = "  +0.30%  "
    "E28" = "  +1.02%  "
    "D29" = "1.01"
    "E29" = "  +0.16%  "
    "D30" = "0.0471"
    "E30" = "  -0.13%  "
    "D31" = "1.11"
    "E31" = "  +0.75%  "
    "E32" = "  -0.18%  "
    "E33" = "  +1.89%  "
    "D34" = "1.421.14"
    "E34" = "  -2.22%  "
    "B35" = "TrustWalletToken"
    "C35" = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
    "D35" = "1.08"
    "E35" = "  +11.83%  "
    "B36" = "LidoDAOToken"
    "C36" = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
    "D36" = "1.61"
    "E36" = "  +2.78%  "
    "D37" = "2.34"
    "E37" = "  +1.73%  "
    "E38" = "  -0.26%  "
    "D39" = "0.532"
    "E39" = "  +1.96%  "
    "B40" = "FraxShare"
    "C40" = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
    "D40" = "5.78"
    "E40" = "  +0.19%  "
    "B41" = "ARBITRUM"
    "C41" = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
    "D41" = "0.809"
    "E41" = "  -0.96%  "
    "E42" = "  +0.01%  "
    "E43" = "  +1.91%  "
    "E44" = "  +1.73%  "
    "D45" = "64.85"
    "E45" = "  +0.30%  "
    "E46" = "  -1.30%  "
    "D47" = "1.698.79"
    "E47" = "  -0.38%  "
    "D48" = "87.56"
    "E48" = "  +0.75%  "
    "B49" = "BabyDogeCoin"
    "C49" = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
    "D49" = "0.0₆0102"
    "E49" = "  +3.05%  "
    "B50" = "Cronos"
    "C50" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
    "D50" = "0.0521"
    "E50" = "  -0.53%  "
    "D51" = "0.0960"
    "E51" = "  -0.70%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
